$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 1100
$ws.Range("B3").Value = "B"

$ws.Range("A3:B3").Select()
